$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.541.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.512.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.511.89'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +3.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.106.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000194'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.518.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.16%  '
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.367.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '421.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.586'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.652.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.34%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000113'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.518.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  -9.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.07%  '
$ws.Range('E38').Value = '  -4.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '173.44'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.36%  '
$ws.Range('E41').Value = '  -5.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.851'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E47').Value = '  -8.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.28%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.902'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.34%  '
